$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the word list (column B, rows 4-15) with the new SRT word set ---
$ws.Range("B4").Value = "THROW"
$ws.Range("B5").Value = "LILY"
$ws.Range("B6").Value = "FILM"
$ws.Range("B7").Value = "DISCREET"
$ws.Range("B8").Value = "LOFT"
$ws.Range("B9").Value = "BEEF"
$ws.Range("B10").Value = "STREET"
$ws.Range("B11").Value = "HELMET"
$ws.Range("B12").Value = "SNAKE"
$ws.Range("B13").Value = "DUG"
$ws.Range("B14").Value = "PACK"
$ws.Range("B15").Value = "TIN"

# These word cells no longer carry the thin box border (matches style used by
# the header cells B1/G1/H1 - font only, no border).
$ws.Range("B4:B15").Borders.LineStyle = -4142

# --- Remove the now-unused "Recog" column (J) entirely ---
$ws.Columns("J").Delete()

# --- Add the new "INTRUSIONS" label row beneath the table ---
# Borrow the plain (font-only, no border) look already used by the header
# cells (B1/G1/H1) rather than the boxed table-cell look.
$ws.Range("B1").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("B17").Value = "INTRUSIONS"
$ws.Rows(17).RowHeight = 45

# --- Restore the selection to match the saved view ---
$ws.Range("J3").Select()

"Done"
